# Apply the "input" sheet data edits + view-state changes described by the
# commit ("bugfixes, it works now!!").
#
# Content changes (sheet "input", row 2):
#   - network_source        (D2): "Network based on shapefile" -> "Network based on OSM online"
#   - OSM_area_of_interest  (E2): (blank)                       -> "zuidholland_4326"
#   - name_of_pbf           (F2): "NL332.osm"                   -> (blank)
#   - network_type          (L2): (blank)                       -> "drive"
#   - road_types            (M2): (blank)                       -> "motorway"
#
# View-state changes:
#   - the "input" sheet becomes the active/selected tab (was "explanation")
#   - selection on "input" moves to M7
#   - selection on "explanation" stays at F1, but it is no longer the active tab

$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("input")

# --- Update the data entry row on the "input" sheet ---
$wsInput.Range("D2").Value = "Network based on OSM online"
$wsInput.Range("E2").Value = "zuidholland_4326"
$wsInput.Range("F2").Value = ""
$wsInput.Range("L2").Value = "drive"
$wsInput.Range("M2").Value = "motorway"

# --- Update view state: make "input" the active sheet/tab and move its selection ---
# (the "explanation" sheet's own selection stays at F1 since we never touch it;
# activating "input" automatically clears tabSelected on the previously active sheet)
$wsInput.Activate()
$wsInput.Range("M7").Select()

# Best-effort window sizing to mirror the author's resized Excel window.
$win = $wb.Windows.Item(1)
$win.Left = -108
$win.Top = -108
$win.Width = 23256
$win.Height = 12576
$win.ScrollColumn = 8
$win.ScrollRow = 1
